$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Valor Mora" total figure
$ws.Range("E11").Value = 1520000

# 2) Update the "Cant. Periodos" count
$ws.Range("F13").Value = 38

# 3) Insert a new detail row before the last (total-bordered) row so that a
#    new period (2509) is added to the statement, while the previously-last
#    row (2508) becomes a normal interior row and the new row takes over the
#    special "closing" border that used to belong to the old last row.
$ws.Rows.Item(52).Insert()

# Copy the formatting (only) of the row above into the newly inserted blank
# row so it matches the other interior detail rows.
$ws.Range("B51:J51").Copy()
$ws.Range("B52:J52").PasteSpecial(-4122)   # xlPasteFormats

# Fill in the data for the duplicated "2508" row (now row 52, interior style)
$ws.Range("B52").Value = "CC"
$ws.Range("C52").Value = "1049026018"
$ws.Range("D52").Value = "DIEGO ANTONIO DIAZ ACEVEDO"
$ws.Range("E52").Value = "2508"
$ws.Range("F52").Value = 40000
$ws.Range("G52").Value = 1000000

# Row 53 already holds the data that used to be in row 52 (it was pushed
# down by the insert, formatting intact). Just update its period to 2509.
$ws.Range("E53").Value = "2509"
$ws.Range("F53").Value = 40000
$ws.Range("G53").Value = 1000000

Write-Output "done"
